$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the "Version" property value (row 3, column B) from 0.4.0 to 0.7.0
$ws.Range("B3").Value = "0.7.0"

# Remove the "Jurisdiction" / "Chile" row entirely (row 11); all rows below shift up by one
$ws.Rows("11").Delete()
